# Auto-generated: apply cryptos.xlsx price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to keep its value as literal text (matches the
    # workbook's inlineStr cells) instead of Excel auto-coercing
    # numeric-looking strings (e.g. "0.9970", "1.140") into numbers,
    # which would silently drop significant trailing/leading zeros.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell 'D2' '27.873.87'
Set-TextCell 'E2' '  +2.77%  '
Set-TextCell 'D3' '1.771.75'
Set-TextCell 'E3' '  -0.54%  '
Set-TextCell 'D4' '0.9995'
Set-TextCell 'E4' '  -0.52%  '
Set-TextCell 'D5' '335.31'
Set-TextCell 'E5' '  -0.50%  '
Set-TextCell 'D6' '0.9970'
Set-TextCell 'E6' '  -0.51%  '
Set-TextCell 'D7' '0.3805'
Set-TextCell 'E7' '  -0.35%  '
Set-TextCell 'D8' '0.3413'
Set-TextCell 'E8' '  +0.15%  '
Set-TextCell 'D9' '47.82'
Set-TextCell 'E9' '  -0.46%  '
Set-TextCell 'D10' '1.140'
Set-TextCell 'E10' '  -3.90%  '
Set-TextCell 'D11' '0.07413'
Set-TextCell 'E11' '  -0.26%  '
Set-TextCell 'B12' 'Solana'
Set-TextCell 'C12' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextCell 'D12' '22.91'
Set-TextCell 'E12' '  +5.93%  '
Set-TextCell 'B13' 'BinanceUSD'
Set-TextCell 'C13' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D13' '0.9967'
Set-TextCell 'E13' '  -0.57%  '
Set-TextCell 'D14' '6.363'
Set-TextCell 'E14' '  -0.91%  '
Set-TextCell 'D15' '1.769.35'
Set-TextCell 'E15' '  -0.49%  '
Set-TextCell 'D16' '7.095'
Set-TextCell 'E16' '  +0.69%  '
Set-TextCell 'D17' '0.00001078'
Set-TextCell 'E17' '  -0.55%  '
Set-TextCell 'D18' '0.06666'
Set-TextCell 'E18' '  +0.37%  '
Set-TextCell 'D19' '82.27'
Set-TextCell 'E19' '  -1.22%  '
Set-TextCell 'D20' '0.9965'
Set-TextCell 'E20' '  -0.57%  '
Set-TextCell 'D21' '17.38'
Set-TextCell 'E21' '  +0.56%  '
Set-TextCell 'D22' '6.408'
Set-TextCell 'E22' '  -2.06%  '
Set-TextCell 'D23' '27.858.34'
Set-TextCell 'D24' '12.07'
Set-TextCell 'E24' '  -1.23%  '
Set-TextCell 'D25' '2.387'
Set-TextCell 'E25' '  +0.25%  '
Set-TextCell 'D26' '1.434'
Set-TextCell 'E26' '  -1.47%  '
Set-TextCell 'D27' '20.70'
Set-TextCell 'E27' '  -1.72%  '
Set-TextCell 'D28' '2.426'
Set-TextCell 'E28' '  -3.23%  '
Set-TextCell 'D29' '153.73'
Set-TextCell 'E29' '  -0.95%  '
Set-TextCell 'D30' '1.969.56'
Set-TextCell 'E30' '  -0.47%  '
Set-TextCell 'D31' '134.23'
Set-TextCell 'E31' '  +0.26%  '
Set-TextCell 'D32' '6.155'
Set-TextCell 'E32' '  +2.29%  '
Set-TextCell 'E33' '  -0.73%  '
Set-TextCell 'D34' '0.08788'
Set-TextCell 'E34' '  +1.41%  '
Set-TextCell 'D35' '12.79'
Set-TextCell 'E35' '  -2.21%  '
Set-TextCell 'D36' '0.02428'
Set-TextCell 'E36' '  +4.92%  '
Set-TextCell 'D37' '0.6871'
Set-TextCell 'E37' '  +0.54%  '
Set-TextCell 'D38' '5.311'
Set-TextCell 'E38' '  -1.24%  '
Set-TextCell 'D39' '0.06331'
Set-TextCell 'E39' '  +0.83%  '
Set-TextCell 'D40' '0.2188'
Set-TextCell 'E40' '  +0.79%  '
Set-TextCell 'D41' '1.517'
Set-TextCell 'E41' '  -6.35%  '
Set-TextCell 'D42' '1.237'
Set-TextCell 'D43' '8.253'
Set-TextCell 'E43' '  -3.60%  '
Set-TextCell 'D44' '14.18'
Set-TextCell 'E44' '  +0.23%  '
Set-TextCell 'D45' '0.9967'
Set-TextCell 'E45' '  -0.53%  '
Set-TextCell 'D46' '0.6290'
Set-TextCell 'E46' '  -2.02%  '
Set-TextCell 'D47' '3.839'
Set-TextCell 'E47' '  -0.42%  '
Set-TextCell 'D48' '131.53'
Set-TextCell 'E48' '  +0.39%  '
Set-TextCell 'D49' '2.085'
Set-TextCell 'D50' '0.07370'
Set-TextCell 'E50' '  +3.91%  '
Set-TextCell 'D51' '1.151'
Set-TextCell 'E51' '  +3.44%  '
